$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.729.69"
$ws.Range("E2").Value = "  +2.47%  "

$ws.Range("D3").Value = "1.864.50"
$ws.Range("E3").Value = "  +2.43%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.038"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +2.09%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "324.38"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.28%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.034"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.48%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4421"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.83%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3797"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.11%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07460"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.74%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8853"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.02%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "21.72"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.21%  "

$ws.Range("D12").Value = "1.874.16"
$ws.Range("E12").Value = "  -15.52%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.553"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.45%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.754"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.97%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.07230"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.52%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "84.00"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.24%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.040"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.25%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000009113"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.92%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.034"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.59%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "15.56"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.91%  "

$ws.Range("D21").Value = "27.745.80"
$ws.Range("E21").Value = "  +2.36%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.315"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.74%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "11.32"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.88%  "

$ws.Range("D24").Value = "2.093.68"
$ws.Range("E24").Value = "  -13.09%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.015"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +6.78%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "158.88"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.82%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.83"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.59%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.988"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.65%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.317"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.79%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "118.14"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.91%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.09056"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.00%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7803"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.51%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.215"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.03%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.022"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +7.56%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.572"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.13%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.037"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.71%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.150"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01994"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.38%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05346"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.47%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.855"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.19%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.5203"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.89%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1691"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.42%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.877"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +6.10%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.667"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.67%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "110.30"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.11%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "10.66"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.71%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.06578"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +5.38%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.721"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.64%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.4714"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.95%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.931"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.64%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "39.82"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.34%  "
